$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object "object[,]" 20,10
$data[0,0] = -17.35316405942421
$data[0,1] = -0.4689816133227405
$data[0,2] = -17.35316405942421
$data[0,3] = -17.35316405942421
$data[0,4] = -17.35316405942421
$data[0,5] = -17.35316405942421
$data[0,6] = -17.35316405942421
$data[0,7] = -17.35316405942421
$data[0,8] = -17.35316405942421
$data[0,9] = -17.35316405942421
$data[1,0] = -17.35316405942421
$data[1,1] = -17.35316405942421
$data[1,2] = -17.35316405942421
$data[1,3] = -17.35316405942421
$data[1,4] = -17.35316405942421
$data[1,5] = -17.35316405942421
$data[1,6] = -17.35316405942421
$data[1,7] = 0.3797519252196256
$data[1,8] = -17.35316405942421
$data[1,9] = -17.35316405942421
$data[2,0] = -17.35316405942421
$data[2,1] = -0.3471338969869318
$data[2,2] = 0.4144320342072859
$data[2,3] = -17.35316405942421
$data[2,4] = 3.899663616282851
$data[2,5] = -17.35316405942421
$data[2,6] = 0.6580031261413791
$data[2,7] = -17.35316405942421
$data[2,8] = 1.62831681961823
$data[2,9] = -17.35316405942421
$data[3,0] = -17.35316405942421
$data[3,1] = 0.2963236965265545
$data[3,2] = -17.35316405942421
$data[3,3] = -17.35316405942421
$data[3,4] = -17.35316405942421
$data[3,5] = 3.440106636019695
$data[3,6] = -17.35316405942421
$data[3,7] = -17.35316405942421
$data[3,8] = -17.35316405942421
$data[3,9] = -17.35316405942421
$data[4,0] = -17.35316405942421
$data[4,1] = -17.35316405942421
$data[4,2] = -17.35316405942421
$data[4,3] = -17.35316405942421
$data[4,4] = -17.35316405942421
$data[4,5] = -17.35316405942421
$data[4,6] = -17.35316405942421
$data[4,7] = -17.35316405942421
$data[4,8] = -17.35316405942421
$data[4,9] = -17.35316405942421
$data[5,0] = 3.024894728596113
$data[5,1] = -17.35316405942421
$data[5,2] = -17.35316405942421
$data[5,3] = -17.35316405942421
$data[5,4] = -17.35316405942421
$data[5,5] = -17.35316405942421
$data[5,6] = -17.35316405942421
$data[5,7] = -17.35316405942421
$data[5,8] = -17.35316405942421
$data[5,9] = -17.35316405942421
$data[6,0] = -17.35316405942421
$data[6,1] = -17.35316405942421
$data[6,2] = -17.35316405942421
$data[6,3] = 2.270154839552673
$data[6,4] = -17.35316405942421
$data[6,5] = -17.35316405942421
$data[6,6] = -17.35316405942421
$data[6,7] = -17.35316405942421
$data[6,8] = -17.35316405942421
$data[6,9] = -17.35316405942421
$data[7,0] = 3.568111035266024
$data[7,1] = -17.35316405942421
$data[7,2] = -17.35316405942421
$data[7,3] = -17.35316405942421
$data[7,4] = -17.35316405942421
$data[7,5] = -17.35316405942421
$data[7,6] = -17.35316405942421
$data[7,7] = -17.35316405942421
$data[7,8] = -17.35316405942421
$data[7,9] = -17.35316405942421
$data[8,0] = -17.35316405942421
$data[8,1] = -17.35316405942421
$data[8,2] = -17.35316405942421
$data[8,3] = -17.35316405942421
$data[8,4] = -17.35316405942421
$data[8,5] = -17.35316405942421
$data[8,6] = -17.35316405942421
$data[8,7] = 0.8604575674235262
$data[8,8] = -17.35316405942421
$data[8,9] = -17.35316405942421
$data[9,0] = -17.35316405942421
$data[9,1] = -17.35316405942421
$data[9,2] = -17.35316405942421
$data[9,3] = 1.99078717313962
$data[9,4] = -17.35316405942421
$data[9,5] = 1.5754612972532
$data[9,6] = -17.35316405942421
$data[9,7] = -17.35316405942421
$data[9,8] = -17.35316405942421
$data[9,9] = -17.35316405942421
$data[10,0] = -17.35316405942421
$data[10,1] = -17.35316405942421
$data[10,2] = -17.35316405942421
$data[10,3] = -17.35316405942421
$data[10,4] = -17.35316405942421
$data[10,5] = -17.35316405942421
$data[10,6] = -17.35316405942421
$data[10,7] = -17.35316405942421
$data[10,8] = -17.35316405942421
$data[10,9] = -17.35316405942421
$data[11,0] = -17.35316405942421
$data[11,1] = -17.35316405942421
$data[11,2] = -17.35316405942421
$data[11,3] = 1.812785651964227
$data[11,4] = -17.35316405942421
$data[11,5] = -17.35316405942421
$data[11,6] = -17.35316405942421
$data[11,7] = -17.35316405942421
$data[11,8] = 1.343331499061982
$data[11,9] = -17.35316405942421
$data[12,0] = -17.35316405942421
$data[12,1] = -17.35316405942421
$data[12,2] = 2.006227338827988
$data[12,3] = -17.35316405942421
$data[12,4] = -17.35316405942421
$data[12,5] = -17.35316405942421
$data[12,6] = -17.35316405942421
$data[12,7] = -17.35316405942421
$data[12,8] = -17.35316405942421
$data[12,9] = -17.35316405942421
$data[13,0] = -17.35316405942421
$data[13,1] = -17.35316405942421
$data[13,2] = -0.5229335969267925
$data[13,3] = -17.35316405942421
$data[13,4] = -17.35316405942421
$data[13,5] = -17.35316405942421
$data[13,6] = -17.35316405942421
$data[13,7] = -17.35316405942421
$data[13,8] = -17.35316405942421
$data[13,9] = -17.35316405942421
$data[14,0] = -17.35316405942421
$data[14,1] = -17.35316405942421
$data[14,2] = -17.35316405942421
$data[14,3] = -17.35316405942421
$data[14,4] = -17.35316405942421
$data[14,5] = -17.35316405942421
$data[14,6] = -17.35316405942421
$data[14,7] = -17.35316405942421
$data[14,8] = 3.16038331945128
$data[14,9] = -17.35316405942421
$data[15,0] = -17.35316405942421
$data[15,1] = 0.1773108625862125
$data[15,2] = -0.1965575321703975
$data[15,3] = -17.35316405942421
$data[15,4] = -17.35316405942421
$data[15,5] = -17.35316405942421
$data[15,6] = 3.113682126390575
$data[15,7] = -0.2439353336347939
$data[15,8] = 1.557481598589201
$data[15,9] = -17.35316405942421
$data[16,0] = -17.35316405942421
$data[16,1] = -17.35316405942421
$data[16,2] = -17.35316405942421
$data[16,3] = -17.35316405942421
$data[16,4] = -17.35316405942421
$data[16,5] = -17.35316405942421
$data[16,6] = 1.482741392409271
$data[16,7] = 0.1974399521978357
$data[16,8] = 1.314380653123583
$data[16,9] = -17.35316405942421
$data[17,0] = -17.35316405942421
$data[17,1] = -17.35316405942421
$data[17,2] = 2.764272323514074
$data[17,3] = -17.35316405942421
$data[17,4] = -17.35316405942421
$data[17,5] = -17.35316405942421
$data[17,6] = 0.6691007931778904
$data[17,7] = 2.948994001899961
$data[17,8] = -17.35316405942421
$data[17,9] = -17.35316405942421
$data[18,0] = -17.35316405942421
$data[18,1] = 3.068975667175011
$data[18,2] = 2.652418003406576
$data[18,3] = -17.35316405942421
$data[18,4] = 2.343363760766668
$data[18,5] = -17.35316405942421
$data[18,6] = 1.552966126643952
$data[18,7] = 2.841969008762656
$data[18,8] = -17.35316405942421
$data[18,9] = 4.321919908798121
$data[19,0] = -17.35316405942421
$data[19,1] = 2.952471341312691
$data[19,2] = -17.35316405942421
$data[19,3] = 2.942686693182929
$data[19,4] = -17.35316405942421
$data[19,5] = 2.624324106386647
$data[19,6] = 1.291047829171505
$data[19,7] = -17.35316405942421
$data[19,8] = -17.35316405942421
$data[19,9] = -17.35316405942421

$ws.Range("B2:K21").Value = $data
